$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status moved from "In Translation" to "Ready for handoff"
$wsOverview.Range("E2:F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Refresh the handoff generation timestamps
$wsOverview.Range("G2").Value = "2016-08-17 08:57:45"
$wsDeDe.Range("H2").Value = "2016-08-17 08:57:45"
$wsZhCn.Range("H2").Value = "2016-08-17 08:57:40"

# Widen the Status columns to fit the new "Ready for handoff" label
# (16.333333333333332 "characters" is the input that Excel's pixel-quantized
# ColumnWidth rounds to the stored width used for these columns)
$wsOverview.Range("E1:F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
